# ---------------------------------------------------------------------------
# Implements: "Launch angle fully implemented; more outputs"
#   - Adds two new single-value inputs to the "Simulation Conditions
#     (Weather)" sheet: "Rail length (effective)" (m) and "Launch angle"
#     (degrees), following the existing input-row pattern/style.
#   - Changes the engine "Expansion ratio" input (Engine Parameters!C6)
#     from 4.49 to 4.99.
#   - Re-points the propellant flow-rate inputs (Propellant Parameters
#     (Tanks)!J6 / J8) from 1593 / 465.5 to 1400 / 400, and switches the
#     "Fuel temperature" row from a Range-of-values entry back to a
#     Single-value entry (290, no step/end).
#   - Leaves every dependent formula cell to recalculate automatically.
#   - Updates sheet selections / the active tab to match the saved UI
#     state after the edit.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsSim     = $wb.Worksheets.Item("Simulation Conditions (Weather)")
$wsRocket  = $wb.Worksheets.Item("Rocket Parameters (Mass)")
$wsEngine  = $wb.Worksheets.Item("Engine Parameters")
$wsTanks   = $wb.Worksheets.Item("Propellant Parameters (Tanks)")

# ---------------------------------------------------------------------------
# 1) Simulation Conditions (Weather): two new input rows (19/20 and 22/23),
#    mirroring the existing "# of Monte Carlo runs" row (16/17) layout:
#    a blank "Value" label row, then the labelled/styled input row.
# ---------------------------------------------------------------------------

# Row 19/20: "Rail length (effective)" = 5 m
$wsSim.Range("C19").Value = "Value"

$wsSim.Range("A17:E17").Copy()
$wsSim.Range("A20:E20").PasteSpecial(-4122)   # xlPasteFormats
$wsSim.Range("A20").Value = "Rail length (effective)"
$wsSim.Range("B20").Value = "Single value"
$wsSim.Range("C20").Value = 5
$wsSim.Range("F20").Value = "m"

# Row 22/23: "Launch angle" = 32 degrees
$wsSim.Range("C22").Value = "Value"

$wsSim.Range("A17:E17").Copy()
$wsSim.Range("A23:E23").PasteSpecial(-4122)   # xlPasteFormats
$wsSim.Range("A23").Value = "Launch angle"
$wsSim.Range("B23").Value = "Single value"
$wsSim.Range("C23").Value = 32
$wsSim.Range("F23").Value = "degrees"

# Same "pick one of the three parameter-mode strings" dropdown as the other
# B-column mode cells (B13 / B10 / B7) on this sheet.
$wsSim.Range("B20").Validation.Add(3, 1, 1, "=Validation!`$A`$2:`$A`$4")
$wsSim.Range("B23").Validation.Add(3, 1, 1, "=Validation!`$A`$2:`$A`$4")

# ---------------------------------------------------------------------------
# 2) Engine Parameters: Expansion ratio 4.49 -> 4.99
# ---------------------------------------------------------------------------

$wsEngine.Range("C6").Value = 4.99

# ---------------------------------------------------------------------------
# 3) Propellant Parameters (Tanks): oxidizer/fuel flow rates, and the Fuel
#    temperature row switching from "Range of values" back to
#    "Single value" (clearing the Step/End columns).
# ---------------------------------------------------------------------------

$wsTanks.Range("J6").Value = 1400
$wsTanks.Range("J8").Value = 400

$wsTanks.Range("B22").Value = "Single value"
$wsTanks.Range("D22").Value = ""
$wsTanks.Range("E22").Value = ""

# ---------------------------------------------------------------------------
# 4) View state: selections on each sheet, and which tab is active.
#    (All downstream formula cells recalc automatically after this script
#    runs, matching the '<v>' changes baked into the target workbook.)
# ---------------------------------------------------------------------------

$wsEngine.Range("C7").Select()
$wsTanks.Range("J9").Select()
$wsSim.Range("C24").Select()
$wsSim.Activate()
